$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 is reference PMC7099829 ("Analysis of clinical characteristics and
# laboratory findings of 95 cases of 2019 novel coronavirus pneumonia in
# Wuhan, China: a retrospective analysis"). Re-running the reference dump
# against this file picked up a cleaned Abstract (the stray
# id="Par1">/id="Par2">/id="Par3">/id="Par4"> markup leftovers from the
# source HTML are now stripped out) and a re-serialized Authors list (extra
# padding spaces around each separator, and Gemin Zhang's author-index
# bumped from 1 to 3) -- fixing the bug from the commit message.
$newAbstract = @"
Background
Since December 2019, 2019 novel coronavirus pneumonia emerged in Wuhan city and rapidly spread throughout China and even the world.
 We sought to analyse the clinical characteristics and laboratory findings of some cases with 2019 novel coronavirus pneumonia .
Methods
In this retrospective study, we extracted the data on 95 patients with laboratory-confirmed 2019 novel coronavirus pneumonia in Wuhan Xinzhou District People’s Hospital from January 16th to February 25th, 2020. Cases were confirmed by real-time RT-PCR and abnormal radiologic findings.
 Outcomes were followed up until March 2th, 2020.
Results
Higher temperature, blood leukocyte count, neutrophil count, neutrophil percentage, C-reactive protein level, D-dimer level, alanine aminotransferase activity, aspartate aminotransferase activity, α - hydroxybutyrate dehydrogenase activity, lactate dehydrogenase activity and creatine kinase activity were related to severe 2019 novel coronavirus pneumonia and composite endpoint, and so were lower lymphocyte count, lymphocyte percentage and total protein level.
 Age below 40 or above 60 years old, male, higher Creatinine level, and lower platelet count also seemed related to severe 2019 novel coronavirus pneumonia and composite endpoint, however the P values were greater than 0.05, which mean under the same condition studies of larger samples are needed in the future.
Conclusion
Multiple factors were related to severe 2019 novel coronavirus pneumonia and composite endpoint, and more related studies are needed in the future.

"@

$newAuthors = "[Gemin%Zhang%NULL%3,     Jie%Zhang%945128911@qq.com%1,     Bowen%Wang%NULL%1,     Xionglin%Zhu%NULL%1,     Qiang%Wang%NULL%1,     Shiming%Qiu%NULL%1]"

$ws.Range("D14").Value2 = $newAbstract
$ws.Range("E14").Value2 = $newAuthors
